$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from serial date 45184 to 45186 for every
# data row (row 1 is the header row; data starts at row 2).
$firstRow = 2
$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Columns S, T, U, V, W, X, Y hold HYPERLINK(...) formulas. Each one gains a
# second "friendly name" argument equal to the row's "Beteckning" (column A)
# value, e.g.:
#   =HYPERLINK("https://.../A 66380-2021.xlsx")
# becomes
#   =HYPERLINK("https://.../A 66380-2021.xlsx", "A 66380-2021")
$linkCols = @(19, 20, 21, 22, 23, 24, 25)  # S..Y

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($beteckning)) {
        continue
    }
    foreach ($c in $linkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        }
    }
}
